# "Analysis and manuscript revision" — reviewer-requested pipeline changes.
#
# The "prev_visits" recoding table gets 6 new rows describing how to
# recode treatment variables collected at prior visits (steroids,
# anti-infectives, macrolides, antiplatelets, anticoagulants and
# immunosuppressants) into boolean yes/no flags. Row 93 already existed
# in the sheet (blank) and simply gets filled in along with the new rows
# 90-92 and 94-95. The "prev_visits" tab also becomes the active sheet
# again (it had been left on "year_visit").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("prev_visits")

# old_var, new_var, trans_fun, args1, args2
$newRows = @(
    @(90, "treat_steroids",     "treat_steroids"),
    @(91, "treat_antiinfec",    "treat_antiinfec"),
    @(92, "macrolides",         "treat_macrolides"),
    @(93, "treat_antiplat",     "treat_antiplat"),
    @(94, "treat_anticoag",     "treat_anticoag"),
    @(95, "treat_immunosuppr",  "treat_immunosuppr")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]        # A: old_var
    $ws1.Cells.Item($r, 2).Value = $row[2]        # B: new_var
    $ws1.Cells.Item($r, 3).Value = "recode_yn"    # C: trans_fun
    $ws1.Cells.Item($r, 4).Value = "F"            # D: args1
    $ws1.Cells.Item($r, 5).Value = "T"            # E: args2
}

# Re-activate "prev_visits" and restore the selection to the last edited
# cell, matching the reviewer's editing session.
$ws1.Activate()
$ws1.Range("B93").Select()
